# Weekly refresh of the "Rabanito" price series (Vega Modelo de Temuco).
# Two new daily price records are prepended to the existing history:
#   - one inserted as the new row 34 (most recent record)
#   - one inserted as the new row 53 (a second, slightly older record)
# Excel's row-insert shifts every existing record below the insertion
# point down by one, which is exactly how this weekly log grows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new row at position 34 (pushes old 34..52 -> 35..53)
$ws.Rows.Item(34).Insert()

$ws.Cells.Item(34, 1).Value = 10
$ws.Cells.Item(34, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(34, 3).Value = "La Araucanía"
$ws.Cells.Item(34, 4).Value = 44567
$ws.Cells.Item(34, 5).Value = 9
$ws.Cells.Item(34, 6).Value = 300000001
$ws.Cells.Item(34, 7).Value = "Rabanito"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 40
$ws.Cells.Item(34, 11).Value = 6000
$ws.Cells.Item(34, 12).Value = 6000
$ws.Cells.Item(34, 13).Value = 6000
$ws.Cells.Item(34, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(34, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(34, 16).Value = 500
$ws.Cells.Item(34, 17).Value = 12
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# --- Insert the second new row at position 53 (pushes old-52's-new-home
#     row 53 -> 54), leaving a fresh blank row 53 for the new record.
$ws.Rows.Item(53).Insert()

$ws.Cells.Item(53, 1).Value = 10
$ws.Cells.Item(53, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(53, 3).Value = "La Araucanía"
$ws.Cells.Item(53, 4).Value = 44568
$ws.Cells.Item(53, 5).Value = 9
$ws.Cells.Item(53, 6).Value = 300000001
$ws.Cells.Item(53, 7).Value = "Rabanito"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 20
$ws.Cells.Item(53, 11).Value = 6000
$ws.Cells.Item(53, 12).Value = 6000
$ws.Cells.Item(53, 13).Value = 6000
$ws.Cells.Item(53, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(53, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(53, 16).Value = 500
$ws.Cells.Item(53, 17).Value = 12
$ws.Cells.Item(53, 18).Value = "Hortaliza"
